$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F ("想去人数")
$updates = @{
    2  = 3424
    3  = 152
    4  = 228
    5  = 1759
    6  = 1668
    8  = 386
    10 = 34
    12 = 43
    13 = 237
    14 = 10
    15 = 52
    16 = 244
    22 = 123
    23 = 29
    24 = 421
    25 = 292
    26 = 123
    27 = 50
    30 = 570
    31 = 2388
    32 = 16
    34 = 493
    35 = 588
    37 = 436
    38 = 243
    39 = 363
    40 = 422
    41 = 555
    42 = 419
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
